$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2124.2
$ws.Range("I98").Value = 2071.3333
$ws.Range("J98").Value = 2600
$ws.Range("K98").Value = 2071.3333
$ws.Range("L98").Value = 2600
$ws.Range("M98").Value = -573.3332999999998
$ws.Range("N98").Value = -5596
$ws.Range("H122").Value = 2124.2
$ws.Range("I122").Value = 2071.3333
$ws.Range("J122").Value = 2600
$ws.Range("K122").Value = 6213.999899999999
$ws.Range("L122").Value = 7800
$ws.Range("M122").Value = -3763.999899999999
$ws.Range("N122").Value = -12700
$ws.Range("H127").Value = 1443.875
$ws.Range("I127").Value = 1367
$ws.Range("J127").Value = 1490
$ws.Range("K127").Value = 4101
$ws.Range("L127").Value = 4470
$ws.Range("M127").Value = 859
$ws.Range("N127").Value = -14390
$ws.Range("H129").Value = 1701.6578
$ws.Range("I129").Value = 586
$ws.Range("J129").Value = 2216.577
$ws.Range("K129").Value = 1758
$ws.Range("L129").Value = 6649.731000000001
$ws.Range("M129").Value = 3242
$ws.Range("N129").Value = -16649.731
$ws.Range("H137").Value = 1283.2576
$ws.Range("I137").Value = 1616.1613
$ws.Range("J137").Value = 988.4
$ws.Range("K137").Value = 4848.4839
$ws.Range("L137").Value = 2965.2
$ws.Range("M137").Value = -2298.4839
$ws.Range("N137").Value = -8065.2
$ws.Range("H138").Value = 1913.2988
$ws.Range("I138").Value = 1125.3773
$ws.Range("J138").Value = 3141.5293
$ws.Range("K138").Value = 3376.1319
$ws.Range("L138").Value = 9424.5879
$ws.Range("M138").Value = 1763.8681
$ws.Range("N138").Value = -19704.5879

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15493.153
$ws.Range("I32").Value = 8905.769
$ws.Range("K32").Value = 8905.769
$ws.Range("M32").Value = -8618.769
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H74").Value = 1543.9375
$ws.Range("I74").Value = 1400.5416
$ws.Range("J74").Value = 1974.125
$ws.Range("K74").Value = 1400.5416
$ws.Range("L74").Value = 1974.125
$ws.Range("M74").Value = -526.5416
$ws.Range("N74").Value = -3722.125
$ws.Range("H77").Value = 1543.9375
$ws.Range("I77").Value = 1400.5416
$ws.Range("J77").Value = 1974.125
$ws.Range("K77").Value = 7002.708000000001
$ws.Range("L77").Value = 9870.625
$ws.Range("M77").Value = -2634.708000000001
$ws.Range("N77").Value = -18606.625
$ws.Range("H132").Value = 2054799.6
$ws.Range("I132").Value = 6715.5
$ws.Range("K132").Value = 20146.5
$ws.Range("M132").Value = -17616.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2813.1333
$ws.Range("I20").Value = 1487.25
$ws.Range("J20").Value = 4328.4287
$ws.Range("K20").Value = 1487.25
$ws.Range("L20").Value = 4328.4287
$ws.Range("M20").Value = -1240.25
$ws.Range("N20").Value = -4822.4287

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 38160800
$ws.Range("J3").Value = 43612056
$ws.Range("L3").Value = 43612056
$ws.Range("N3").Value = -43612282
$ws.Range("H31").Value = 1624.1515
$ws.Range("I31").Value = 1034.3286
$ws.Range("J31").Value = 3047.862
$ws.Range("K31").Value = 1034.3286
$ws.Range("L31").Value = 3047.862
$ws.Range("M31").Value = -739.3286
$ws.Range("N31").Value = -3637.862
$ws.Range("H34").Value = 1624.1515
$ws.Range("I34").Value = 1034.3286
$ws.Range("J34").Value = 3047.862
$ws.Range("K34").Value = 1034.3286
$ws.Range("L34").Value = 3047.862
$ws.Range("M34").Value = -832.3286
$ws.Range("N34").Value = -3451.862
$ws.Range("H58").Value = 8523.235
$ws.Range("I58").Value = 4949.5
$ws.Range("J58").Value = 13628.571
$ws.Range("K58").Value = 4949.5
$ws.Range("L58").Value = 13628.571
$ws.Range("M58").Value = -4746.5
$ws.Range("N58").Value = -14034.571
$ws.Range("H62").Value = 10732.083
$ws.Range("I62").Value = 2495.8333
$ws.Range("J62").Value = 18968.334
$ws.Range("K62").Value = 2495.8333
$ws.Range("L62").Value = 18968.334
$ws.Range("M62").Value = -1871.8333
$ws.Range("N62").Value = -20216.334
$ws.Range("H65").Value = 10732.083
$ws.Range("I65").Value = 2495.8333
$ws.Range("J65").Value = 18968.334
$ws.Range("K65").Value = 12479.1665
$ws.Range("L65").Value = 94841.67
$ws.Range("M65").Value = -9359.1665
$ws.Range("N65").Value = -101081.67
$ws.Range("H134").Value = 3667.1428
$ws.Range("I134").Value = 2935
$ws.Range("J134").Value = 3789.1667
$ws.Range("K134").Value = 8805
$ws.Range("L134").Value = 11367.5001
$ws.Range("M134").Value = -6270
$ws.Range("N134").Value = -16437.5001
$ws.Range("H136").Value = 8523.235
$ws.Range("I136").Value = 4949.5
$ws.Range("J136").Value = 13628.571
$ws.Range("K136").Value = 14848.5
$ws.Range("L136").Value = 40885.713
$ws.Range("M136").Value = -12298.5
$ws.Range("N136").Value = -45985.713

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 6249.5
$ws.Range("I109").Value = 2999
$ws.Range("J109").Value = 9500
$ws.Range("K109").Value = 8997
$ws.Range("L109").Value = 28500
$ws.Range("M109").Value = -7957
$ws.Range("N109").Value = -30580
$ws.Range("H127").Value = 8265389
$ws.Range("I127").Value = 700
$ws.Range("J127").Value = 9091858
$ws.Range("K127").Value = 2100
$ws.Range("L127").Value = 27275574
$ws.Range("M127").Value = 2860
$ws.Range("N127").Value = -27285494
$ws.Range("H132").Value = 1134.95
$ws.Range("I132").Value = 1013.06665
$ws.Range("J132").Value = 1208.08
$ws.Range("K132").Value = 9117.59985
$ws.Range("L132").Value = 10872.72
$ws.Range("M132").Value = -6587.599850000001
$ws.Range("N132").Value = -15932.72

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H113").Value = 1760.56
$ws.Range("I113").Value = 1652.7826
$ws.Range("K113").Value = 1652.7826
$ws.Range("M113").Value = 517.2174
$ws.Range("H122").Value = 3291.5715
$ws.Range("I122").Value = 3171.0908
$ws.Range("K122").Value = 9513.2724
$ws.Range("M122").Value = -7063.2724
$ws.Range("H126").Value = 5959298
$ws.Range("I126").Value = 11370010
$ws.Range("J126").Value = 7514.7
$ws.Range("K126").Value = 34110030
$ws.Range("L126").Value = 22544.1
$ws.Range("M126").Value = -34107560
$ws.Range("N126").Value = -27484.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1887.8096
$ws.Range("I82").Value = 1321.909
$ws.Range("J82").Value = 2510.3
$ws.Range("K82").Value = 1321.909
$ws.Range("L82").Value = 2510.3
$ws.Range("M82").Value = -960.9090000000001
$ws.Range("N82").Value = -3232.3
$ws.Range("H85").Value = 1887.8096
$ws.Range("I85").Value = 1321.909
$ws.Range("J85").Value = 2510.3
$ws.Range("K85").Value = 1321.909
$ws.Range("L85").Value = 2510.3
$ws.Range("M85").Value = -73.9090000000001
$ws.Range("N85").Value = -5006.3
$ws.Range("H94").Value = 18436.25
$ws.Range("J94").Value = 18436.25
$ws.Range("L94").Value = 18436.25
$ws.Range("N94").Value = -19788.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 60172.332
$ws.Range("I15").Value = 200006
$ws.Range("J15").Value = 32205.6
$ws.Range("K15").Value = 200006
$ws.Range("L15").Value = 32205.6
$ws.Range("M15").Value = -199718
$ws.Range("N15").Value = -32781.6
$ws.Range("H54").Value = 8000
$ws.Range("J54").Value = 8000
$ws.Range("L54").Value = 8000
$ws.Range("N54").Value = -9040
$ws.Range("H81").Value = 2035.4706
$ws.Range("J81").Value = 2479.8
$ws.Range("L81").Value = 4959.6
$ws.Range("N81").Value = -7081.6
$ws.Range("H84").Value = 2035.4706
$ws.Range("J84").Value = 2479.8
$ws.Range("L84").Value = 24798
$ws.Range("N84").Value = -35406
$ws.Range("H132").Value = 1953.0625
$ws.Range("I132").Value = 1640.1613
$ws.Range("J132").Value = 2523.647
$ws.Range("K132").Value = 4920.4839
$ws.Range("L132").Value = 7570.941
$ws.Range("M132").Value = -2390.4839
$ws.Range("N132").Value = -12630.941
